$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as plain text in the source data (e.g.
# "62.870.19" uses a dot as a thousands separator), so force text formatting
# before assigning to stop Excel from reinterpreting them as numbers, then
# restore the default "Normal" style so no stray formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.870.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.042.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.94%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.041.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("E13").Value = "  -3.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.50%  "
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.549.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.876.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.045.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.01%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0802"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "422.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.283"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.825.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("E51").Value = "  -0.84%  "
